$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain plain text so values like
# "22.464.06" or "0.00001123" are not reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.464.06"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.57"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.16"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3740"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.94"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3405"
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07559"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.992"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.945"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.572.21"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001123"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.95"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06731"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.43"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.18"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.453.85"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.363"
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.587"
$ws.Range("E26").Value = "  -6.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.14"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.75"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.019"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.83"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.747.08"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.141"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.976"
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.842"
$ws.Range("E35").Value = "  -2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08412"
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02464"
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2297"
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06545"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.484"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.35"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6278"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.04"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5846"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.225"
$ws.Range("E50").Value = "  -5.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07326"
$ws.Range("E51").Value = "  -0.14%  "

# Rows 48 and 49 swapped contents (Quant <-> NEARProtocol)
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.088"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.70"
$ws.Range("E49").Value = "  +3.12%  "
